$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new "2020" column (K), mirroring the formatting of
# the existing "2019" column (J) for each affected row.
$ws.Range("J3:J4").Copy()
$ws.Range("K3:K4").PasteSpecial(-4122)

$ws.Range("J6:J8").Copy()
$ws.Range("K6:K8").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# New 2020 values
$ws.Range("K4").Value = 2020
$ws.Range("K6").Value = 5.9
$ws.Range("K7").Value = 1.5
$ws.Range("K8").Value = "-"

# Match the cursor position left after the edit
$ws.Range("L16").Select()
